# Commit: fetch data from DB and generate excel file
#
# - Add a new "__v" column (H) with value 0 for every data row.
# - Reformat the createdDate (B) / dueDate (C) text columns from
#   "M/D/YYYY" to an ISO-8601 UTC instant string, shifting the date back
#   one day (mirrors a JS `new Date(dateString).toISOString()` round-trip
#   through a UTC+5:30 locale) e.g. "3/27/2018" -> "2018-03-26T18:30:00.000Z".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 101

# --- New "__v" column -------------------------------------------------
$ws.Range("H1").Value = "__v"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}

# --- Reformat createdDate (B) / dueDate (C) ----------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @("B", "C")) {
        $cell = $ws.Range("$col$r")
        $raw = [string]$cell.Value2
        $parts = $raw.Split("/")
        $month = [int]$parts[0]
        $day = [int]$parts[1]
        $year = [int]$parts[2]

        $dt = Get-Date -Year $year -Month $month -Day $day -Hour 0 -Minute 0 -Second 0
        $dt = $dt.AddDays(-1)

        $iso = $dt.ToString("yyyy-MM-ddT18:30:00.000Z")
        $cell.Value = $iso
    }
}
